$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$sh = $s.Shapes.Item(2)

$tr0 = $sh.TextFrame.TextRange
$tr0.InsertBefore("`r`r`r`r")
$final = $sh.TextFrame.TextRange

$texts = @(
    "Requisitos del programa",
    "Funcionalidades",
    "Estructura del programa (división de paquetes y clases)",
    "Diagrama de clases"
)
for ($i = 1; $i -le 4; $i++) {
    $cur = $sh.TextFrame.TextRange
    $para = $cur.Paragraphs($i, 1)
    $para.Text = $texts[$i - 1]
}

$finalAll = $sh.TextFrame.TextRange
Write-Host "final: [$($finalAll.Text)] len=$($finalAll.Length)"
